$d = $word.ActiveDocument

# Locate the two target paragraphs by their distinctive text.
$paraCount = $d.Paragraphs.Count
$decorationParaIndex = -1
$sliderBulletParaIndex = -1
for ($i = 1; $i -le $paraCount; $i++) {
    $ptext = $d.Paragraphs.Item($i).Range.Text
    if ($ptext -match "move color into") {
        $decorationParaIndex = $i
    }
    if ($ptext -match "^Slider:") {
        $sliderBulletParaIndex = $i + 1
    }
}

# --- 1) Move the "_GoBack" bookmark to the end of the "...move color into
#        decoration." paragraph (right before its paragraph mark), by
#        replacing that whole paragraph with an identical copy plus the
#        bookmark appended at the end. ---
$pDecoration = $d.Paragraphs.Item($decorationParaIndex)
$rDecoration = $pDecoration.Range

$decorationXml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + `
'<w:p>' + `
  '<w:pPr>' + `
    '<w:pStyle w:val="ListParagraph"/>' + `
    '<w:numPr><w:ilvl w:val="3"/><w:numId w:val="2"/></w:numPr>' + `
    '<w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr>' + `
  '</w:pPr>' + `
  '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve">If have </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/>' + `
  '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>boxDecoration</w:t></w:r>' + `
  '<w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve"> we have to move color into &#8220;decoration&#8221;.</w:t></w:r>' + `
  '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' + `
'</w:p>' + `
'</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$rDecoration.InsertXML($decorationXml)

# --- 2) Re-resolve paragraph indices (paragraph count unchanged by the
#        replacement above) and fill in the previously-empty Slider
#        sub-bullet paragraph (which used to only host the bookmark) with
#        the new "Custom by SliderTheme(widgets)" content. ---
$pSlider = $d.Paragraphs.Item($sliderBulletParaIndex)
$rSlider = $pSlider.Range

$sliderXml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + `
'<w:p>' + `
  '<w:pPr>' + `
    '<w:pStyle w:val="ListParagraph"/>' + `
    '<w:numPr><w:ilvl w:val="3"/><w:numId w:val="2"/></w:numPr>' + `
    '<w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="vi-VN"/></w:rPr>' + `
  '</w:pPr>' + `
  '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve">Custom by </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/>' + `
  '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>SliderTheme</w:t></w:r>' + `
  '<w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>(widgets)</w:t></w:r>' + `
'</w:p>' + `
'</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$rSlider.InsertXML($sliderXml)
